$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    "B2" = -0.185198527466267
    "C2" = 0.6817654794757074
    "D2" = 0.8357345580027106
    "E2" = 0.9141851880241282
    "F2" = 0.9153496465291245
    "G2" = 23

    "B3" = 0.6061131440253159
    "C3" = 0.9235219194524696
    "D3" = 2.046785243796049
    "E3" = 1.430659024294765
    "F3" = 1.32641777175159
    "G3" = 22

    "B4" = 0.4736973475431321
    "C4" = 1.172664837566667
    "D4" = 3.405130736559959
    "E4" = 1.845299633273675
    "F4" = 1.827506125151183
    "G4" = 21

    "B5" = 0.5953765531118547
    "C5" = 0.8086914423054233
    "D5" = 1.001718343406927
    "E5" = 1.000858802932225
    "F5" = 0.8254154318705684
    "G5" = 20

    "B6" = 0.4684946891824553
    "C6" = 0.6936936078512286
    "D6" = 0.7365785875710174
    "E6" = 0.8582415671423852
    "F6" = 0.7387953769405898
    "G6" = 19

    "B7" = 0.2823731533649205
    "C7" = 0.5667961785861276
    "D7" = 0.457481639405574
    "E7" = 0.6763738902453095
    "F7" = 0.6324298031243258
    "G7" = 18

    "B8" = 0.2744592162286795
    "C8" = 0.5038470190539809
    "D8" = 0.354955400894988
    "E8" = 0.5957813364775604
    "F8" = 0.5450727114268883
    "G8" = 17

    "B9" = 0.2856656093728968
    "C9" = 0.4249183470701322
    "D9" = 0.2447576538564426
    "E9" = 0.4947298796883433
    "F9" = 0.4218825517020163
    "G9" = 12

    "B10" = 0.1512223912673487
    "C10" = 0.384830859019198
    "D10" = 0.2185051035033979
    "E10" = 0.467445294663876
    "F10" = 0.4777478838570751
    "G10" = 7

    "B11" = 0.1070798832976198
    "C11" = 0.5676804838982221
    "D11" = 0.3539232406442811
    "E11" = 0.5949144817906865
    "F11" = 0.7167187097152367
    "G11" = 3
}

foreach ($addr in $data.Keys) {
    $ws.Range($addr).Value = $data[$addr]
}
